$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATMI ligand-receptor scoring results re-run following Dr Hou's advice.
# The number of ligand/receptor expressing cells for this Apoe-Lrp5 pair
# (columns E "Ligand-expressing cells" and K "Receptor-expressing cells")
# increased from 1 to 3 for every row, which cascades into the expression,
# specificity and edge-weight statistics in columns G,H,I,J,M,N,O,P,Q,R,S,T.
# Columns F and L (detection rates) are unaffected.
$rowUpdates = @{
    2 = @{ "E"=3; "G"=156.656447; "H"=469.969341; "I"=0.0671576211124673; "J"=0.0671576211124673; "K"=3; "M"=13.523597; "N"=40.570791; "O"=0.3454494697445509; "P"=0.3454494697445509; "Q"=2118.558656679859; "R"=19067.02791011873; "S"=0.02319956460260729; "T"=0.02319956460260729 }
    3 = @{ "E"=3; "G"=156.656447; "H"=469.969341; "I"=0.0671576211124673; "J"=0.0671576211124673; "K"=3; "M"=11.59690133333333; "N"=34.790704; "O"=0.2962335698320407; "P"=0.2962335698320406; "Q"=1816.729359089562; "R"=16350.56423180606; "S"=0.01989434184357381; "T"=0.01989434184357381 }
    4 = @{ "E"=3; "G"=156.656447; "H"=469.969341; "I"=0.0671576211124673; "J"=0.0671576211124673; "K"=3; "M"=4.574811666666666; "N"=13.724435; "O"=0.1168599052775075; "P"=0.1168599052775075; "Q"=716.6737413941482; "R"=6450.063672547335; "S"=0.007848033241865669; "T"=0.007848033241865669 }
    5 = @{ "E"=3; "G"=156.656447; "H"=469.969341; "I"=0.0671576211124673; "J"=0.0671576211124673; "K"=3; "M"=9.452519666666667; "N"=28.357559; "O"=0.2414570551459009; "P"=0.2414570551459009; "Q"=1480.798146177624; "R"=13327.18331559862; "S"=0.01621568142442053; "T"=0.01621568142442053 }
    6 = @{ "E"=3; "G"=42.300692; "H"=126.902076; "I"=0.01813403725498241; "J"=0.01813403725498241; "K"=3; "M"=13.523597; "N"=40.570791; "O"=0.3454494697445509; "P"=0.3454494697445509; "Q"=572.057511429124; "R"=5148.517602862115; "S"=0.006264393554061604; "T"=0.006264393554061604 }
    7 = @{ "E"=3; "G"=42.300692; "H"=126.902076; "I"=0.01813403725498241; "J"=0.01813403725498241; "K"=3; "M"=11.59690133333333; "N"=34.790704; "O"=0.2962335698320407; "P"=0.2962335698320406; "Q"=490.5569514557226; "R"=4415.012563101503; "S"=0.005371910591510658; "T"=0.005371910591510657 }
    8 = @{ "E"=3; "G"=42.300692; "H"=126.902076; "I"=0.01813403725498241; "J"=0.01813403725498241; "K"=3; "M"=4.574811666666666; "N"=13.724435; "O"=0.1168599052775075; "P"=0.1168599052775075; "Q"=193.5176992696733; "R"=1741.65929342706; "S"=0.002119141875916037; "T"=0.002119141875916037 }
    9 = @{ "E"=3; "G"=42.300692; "H"=126.902076; "I"=0.01813403725498241; "J"=0.01813403725498241; "K"=3; "M"=9.452519666666667; "N"=28.357559; "O"=0.2414570551459009; "P"=0.2414570551459009; "Q"=399.8481230436093; "R"=3598.633107392484; "S"=0.004378591233494108; "T"=0.004378591233494108 }
    10 = @{ "E"=3; "G"=2110.189616; "H"=6330.568848; "I"=0.9046248489651427; "J"=0.9046248489651426; "K"=3; "M"=13.523597; "N"=40.570791; "O"=0.3454494697445509; "P"=0.3454494697445509; "Q"=28537.35396036875; "R"=256836.1856433187; "S"=0.312502174392753; "T"=0.312502174392753 }
    11 = @{ "E"=3; "G"=2110.189616; "H"=6330.568848; "I"=0.9046248489651427; "J"=0.9046248489651426; "K"=3; "M"=11.59690133333333; "N"=34.790704; "O"=0.2962335698320407; "P"=0.2962335698320406; "Q"=24471.66077137656; "R"=220244.946942389; "S"=0.2679802483677149; "T"=0.2679802483677148 }
    12 = @{ "E"=3; "G"=2110.189616; "H"=6330.568848; "I"=0.9046248489651427; "J"=0.9046248489651426; "K"=3; "M"=4.574811666666666; "N"=13.724435; "O"=0.1168599052775075; "P"=0.1168599052775075; "Q"=9653.720074155653; "R"=86883.48066740087; "S"=0.1057143741617461; "T"=0.1057143741617461 }
    13 = @{ "E"=3; "G"=2110.189616; "H"=6330.568848; "I"=0.9046248489651427; "J"=0.9046248489651426; "K"=3; "M"=9.452519666666667; "N"=28.357559; "O"=0.2414570551459009; "P"=0.2414570551459009; "Q"=19946.60884563578; "R"=179519.479610722; "S"=0.2184280520429287; "T"=0.2184280520429287 }
    14 = @{ "E"=3; "G"=23.52144266666667; "H"=70.564328; "I"=0.01008349266740757; "J"=0.01008349266740757; "K"=3; "M"=13.523597; "N"=40.570791; "O"=0.3454494697445509; "P"=0.3454494697445509; "Q"=318.0945114826054; "R"=2862.850603343448; "S"=0.003483337195129013; "T"=0.003483337195129012 }
    15 = @{ "E"=3; "G"=23.52144266666667; "H"=70.564328; "I"=0.01008349266740757; "J"=0.01008349266740757; "K"=3; "M"=11.59690133333333; "N"=34.790704; "O"=0.2962335698320407; "P"=0.2962335698320406; "Q"=272.7758498229903; "R"=2454.982648406912; "S"=0.002987069029241351; "T"=0.00298706902924135 }
    16 = @{ "E"=3; "G"=23.52144266666667; "H"=70.564328; "I"=0.01008349266740757; "J"=0.01008349266740757; "K"=3; "M"=4.574811666666666; "N"=13.724435; "O"=0.1168599052775075; "P"=0.1168599052775075; "Q"=107.6061703282978; "R"=968.4555329546801; "S"=0.001178355997979691; "T"=0.00117835599797969 }
    17 = @{ "E"=3; "G"=23.52144266666667; "H"=70.564328; "I"=0.01008349266740757; "J"=0.01008349266740757; "K"=3; "M"=9.452519666666667; "N"=28.357559; "O"=0.2414570551459009; "P"=0.2414570551459009; "Q"=222.3368993950392; "R"=2001.032094555352; "S"=0.002434730445057517; "T"=0.002434730445057517 }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}